$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 8, pushing old rows 8-13
# down to rows 10-15.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# New row 8: Terminal Hortofrutícola Agro Chillán - Perejil, Primera, 2022-08-24
$ws.Cells.Item(8,1).Value  = 7
$ws.Cells.Item(8,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8,3).Value  = "Ñuble"
$ws.Cells.Item(8,4).Value  = 44797
$ws.Cells.Item(8,5).Value  = 16
$ws.Cells.Item(8,6).Value  = 100112044
$ws.Cells.Item(8,7).Value  = "Perejil"
$ws.Cells.Item(8,8).Value  = "Sin especificar"
$ws.Cells.Item(8,9).Value  = "Primera"
$ws.Cells.Item(8,10).Value = 240
$ws.Cells.Item(8,11).Value = 750
$ws.Cells.Item(8,12).Value = 850
$ws.Cells.Item(8,13).Value = 800
$ws.Cells.Item(8,14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(8,15).Value = "Región del Maule"
$ws.Cells.Item(8,16).Value = 800
$ws.Cells.Item(8,17).Value = 1
$ws.Cells.Item(8,18).Value = "Hortaliza"

# New row 9: Terminal Hortofrutícola Agro Chillán - Perejil, Segunda, 2022-08-24
$ws.Cells.Item(9,1).Value  = 7
$ws.Cells.Item(9,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9,3).Value  = "Ñuble"
$ws.Cells.Item(9,4).Value  = 44797
$ws.Cells.Item(9,5).Value  = 16
$ws.Cells.Item(9,6).Value  = 100112044
$ws.Cells.Item(9,7).Value  = "Perejil"
$ws.Cells.Item(9,8).Value  = "Sin especificar"
$ws.Cells.Item(9,9).Value  = "Segunda"
$ws.Cells.Item(9,10).Value = 200
$ws.Cells.Item(9,11).Value = 650
$ws.Cells.Item(9,12).Value = 650
$ws.Cells.Item(9,13).Value = 650
$ws.Cells.Item(9,14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(9,15).Value = "Región del Maule"
$ws.Cells.Item(9,16).Value = 650
$ws.Cells.Item(9,17).Value = 1
$ws.Cells.Item(9,18).Value = "Hortaliza"

# Apply the same date-number formatting used by column D elsewhere.
$ws.Cells.Item(8,4).NumberFormat = $ws.Cells.Item(11,4).NumberFormat
$ws.Cells.Item(9,4).NumberFormat = $ws.Cells.Item(11,4).NumberFormat
